$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume(1h) columns store scraped values as plain text (e.g.
# "40.31", "5.99%"). Forcing the Text number format before writing keeps
# Excel from auto-coercing these into numeric/percentage values, which would
# change the underlying cell type and the literal text Excel reports back.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "328.33"
$ws.Range("E2").Value = "5.99%"
$ws.Range("D3").Value = "40.31"
$ws.Range("E3").Value = "8.84%"
$ws.Range("D4").Value = "5.601"
$ws.Range("E4").Value = "9.29%"
$ws.Range("D5").Value = "0.08161"
$ws.Range("E5").Value = "3.97%"
$ws.Range("D6").Value = "4.547"
$ws.Range("E6").Value = "3.40%"
$ws.Range("D7").Value = "8.671"
$ws.Range("E7").Value = "4.84%"
$ws.Range("D8").Value = "1.975"
$ws.Range("E8").Value = "5.43%"
$ws.Range("E9").Value = "-0.15%"
$ws.Range("D10").Value = "0.9490"
$ws.Range("E10").Value = "2.73%"
$ws.Range("E11").Value = "9.62%"
$ws.Range("D12").Value = "0.1983"
$ws.Range("E12").Value = "4.55%"
$ws.Range("D13").Value = "0.09167"
$ws.Range("E13").Value = "2.87%"
$ws.Range("D14").Value = "0.03568"
$ws.Range("E14").Value = "7.54%"
$ws.Range("D15").Value = "0.09610"
$ws.Range("E15").Value = "0.17%"
$ws.Range("D16").Value = "0.001309"
$ws.Range("E16").Value = "-5.47%"
$ws.Range("D17").Value = "0.006223"
$ws.Range("E17").Value = "0.48%"
$ws.Range("D18").Value = "3.370"
$ws.Range("E18").Value = "-0.61%"
$ws.Range("D19").Value = "0.3519"
$ws.Range("E19").Value = "1.76%"
$ws.Range("D20").Value = "7.479"
$ws.Range("E20").Value = "17.42%"
$ws.Range("D21").Value = "0.1397"
$ws.Range("E21").Value = "7.67%"
$ws.Range("D22").Value = "0.2481"
$ws.Range("E22").Value = "3.27%"
$ws.Range("D23").Value = "0.04426"
$ws.Range("E23").Value = "1.93%"
$ws.Range("D24").Value = "0.001257"
$ws.Range("E24").Value = "4.62%"
$ws.Range("D25").Value = "0.004293"
$ws.Range("E25").Value = "0.17%"
$ws.Range("E26").Value = "-15.20%"
$ws.Range("D27").Value = "0.0003983"
$ws.Range("E27").Value = "37.41%"
$ws.Range("D39").Value = "0.02539"
$ws.Range("E39").Value = "17.27%"
$ws.Range("D40").Value = "0.05214"
$ws.Range("E40").Value = "4.12%"
$ws.Range("D41").Value = "0.007788"
$ws.Range("E41").Value = "2.60%"
$ws.Range("D42").Value = "0.1437"
$ws.Range("E42").Value = "6.05%"
$ws.Range("D43").Value = "0.008988"
$ws.Range("E43").Value = "5.49%"
$ws.Range("D44").Value = "0.002187"
$ws.Range("E44").Value = "8.83%"
$ws.Range("D45").Value = "0.009605"
$ws.Range("E45").Value = "7.77%"
$ws.Range("D46").Value = "0.00006707"
$ws.Range("E46").Value = "2.13%"
$ws.Range("E47").Value = "-0.26%"
$ws.Range("D48").Value = "0.002869"
$ws.Range("E48").Value = "-12.96%"
$ws.Range("E49").Value = "59.43%"
$ws.Range("E50").Value = "-0.26%"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").Value = "-0.26%"
